$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    # Force the cell to remain a literal text string (the source workbook
    # stores these as inline strings) instead of letting Excel's COM layer
    # auto-coerce percent-looking / date-looking / numeric-looking text
    # into a number, date or time serial.
    $range.NumberFormat = "@"
    $range.Value = $value
}

# ---------------------------------------------------------------------------
# 1. Summary sheet - update OVERALL and leadlag strategy rows with new totals
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("C2").Value = 7
Set-TextValue $summary.Range("D2") "57.1%"
Set-TextValue $summary.Range("E2") "+1.3828%"
Set-TextValue $summary.Range("F2") "+0.1975%"

$summary.Range("C3").Value = 19
Set-TextValue $summary.Range("D3") "15.8%"
Set-TextValue $summary.Range("E3") "+1.3375%"
Set-TextValue $summary.Range("F3") "+0.0704%"

# ---------------------------------------------------------------------------
# 2. leadlag sheet - close trade #7 (row 6) and append new trade #25 (row 21)
# ---------------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

$leadlag.Range("G6").Value = 69840.50732400001
Set-TextValue $leadlag.Range("H6") "CLOSED"
$leadlag.Range("I6").Value = 0.5741000000000001
$leadlag.Range("J6").Value = 5.74
Set-TextValue $leadlag.Range("M6") "time_exit_5min"
$leadlag.Range("N6").Value = 5

$leadlag.Range("A21").Value = 25
Set-TextValue $leadlag.Range("B21") "2026-02-16"
Set-TextValue $leadlag.Range("C21") "21:26:26"
Set-TextValue $leadlag.Range("D21") "leadlag"
Set-TextValue $leadlag.Range("E21") "DOWN"
$leadlag.Range("F21").Value = 69021.13
Set-TextValue $leadlag.Range("H21") "OPEN"
$leadlag.Range("I21").Value = 0
$leadlag.Range("J21").Value = 0
$leadlag.Range("K21").Value = 0.6344
Set-TextValue $leadlag.Range("L21") "Binance leading with -0.063% move"
$leadlag.Range("N21").Value = 0

# ---------------------------------------------------------------------------
# 3. All Trades sheet - append closed trade #7 as new row 8
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("A8").Value = 7
Set-TextValue $allTrades.Range("B8") "2026-02-16"
Set-TextValue $allTrades.Range("C8") "21:21:22"
Set-TextValue $allTrades.Range("D8") "leadlag"
Set-TextValue $allTrades.Range("E8") "UP"
$allTrades.Range("F8").Value = 69441.86
$allTrades.Range("G8").Value = 69840.50732400001
Set-TextValue $allTrades.Range("H8") "CLOSED"
$allTrades.Range("I8").Value = 0.5741000000000001
$allTrades.Range("J8").Value = 5.74
$allTrades.Range("K8").Value = 0.75
Set-TextValue $allTrades.Range("L8") "Binance leading with 0.079% move"
Set-TextValue $allTrades.Range("M8") "time_exit_5min"
$allTrades.Range("N8").Value = 5

# ---------------------------------------------------------------------------
# 4. Comparison sheet - refresh leadlag aggregate stats
# ---------------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

$comparison.Range("B2").Value = 19
Set-TextValue $comparison.Range("C2") "15.8%"
Set-TextValue $comparison.Range("D2") "3.82"
Set-TextValue $comparison.Range("E2") "+0.6038%"
Set-TextValue $comparison.Range("G2") "2.55"
